# Weekly update: insert the latest week's price record for Ají (Inferno, Primera)
# at the top of the "Feria Lagunitas de Puerto Montt" price table, pushing all
# subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 65, shifting rows 65:158 down to 66:159.
$ws.Rows(65).Insert()

# Populate the new row with the new weekly record.
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44495
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = 100112021
$ws.Range("G65").Value = "Ají"
$ws.Range("H65").Value = "Inferno"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 160
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = 35000
$ws.Range("N65").Value = "$/caja 12 kilos"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 2917
$ws.Range("Q65").Value = 12
$ws.Range("R65").Value = "Hortaliza"
